# Generate Report for Handoff
# Adds a new localization-status row for e1985793-3bf4-49e6-81c5-3825c0ba1c49
# to the Overview, zh-cn and de-de worksheets.

$wb = $excel.ActiveWorkbook

$fileGuid = "e1985793-3bf4-49e6-81c5-3825c0ba1c49"
$mdName = "$fileGuid.md"
$zhXlf  = "$fileGuid.c115b38d8653e4933062376f45500448e1bda95c.zh-cn.xlf"
$deXlf  = "$fileGuid.c115b38d8653e4933062376f45500448e1bda95c.de-de.xlf"

$status        = "Ready for handoff"
$ext           = ".md"
$reason        = "Include"
$noHandback    = "0001-01-01 00:00:00"
$handoffDate   = "2016-03-21 10:38:49"
$zhHandoffTime = "2016-03-21 10:38:46"
$deHandoffTime = "2016-03-21 10:38:49"

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A9").Value = $mdName
$wsOverview.Range("B9").Value = $status
$wsOverview.Range("C9").Value = $status
$wsOverview.Range("D9").Value = $handoffDate
$wsOverview.Range("D9").NumberFormat = $dateFmt

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A9"),
    "https://github.com/OpenLocalizationTest/oltest/blob/0000000000000000000000000000000000000000/e2e/$mdName",
    "",
    "",
    $mdName
)

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A9").Value = $mdName
$wsZh.Range("B9").Value = $ext
$wsZh.Range("C9").Value = $status
$wsZh.Range("D9").Value = $zhXlf
$wsZh.Range("E9").Value = $zhHandoffTime
$wsZh.Range("H9").Value = $noHandback
$wsZh.Range("J9").Value = $reason

$wsZh.Range("E9").NumberFormat = $dateFmt
$wsZh.Range("H9").NumberFormat = $dateFmt

$wsZh.Hyperlinks.Add(
    $wsZh.Range("A9"),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/0000000000000000000000000000000000000000/e2e/$mdName",
    "",
    "",
    $mdName
)
$wsZh.Hyperlinks.Add(
    $wsZh.Range("D9"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlf",
    "",
    "",
    $zhXlf
)

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A9").Value = $mdName
$wsDe.Range("B9").Value = $ext
$wsDe.Range("C9").Value = $status
$wsDe.Range("D9").Value = $deXlf
$wsDe.Range("E9").Value = $deHandoffTime
$wsDe.Range("H9").Value = $noHandback
$wsDe.Range("J9").Value = $reason

$wsDe.Range("E9").NumberFormat = $dateFmt
$wsDe.Range("H9").NumberFormat = $dateFmt

$wsDe.Hyperlinks.Add(
    $wsDe.Range("A9"),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/0000000000000000000000000000000000000000/e2e/$mdName",
    "",
    "",
    $mdName
)
$wsDe.Hyperlinks.Add(
    $wsDe.Range("D9"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlf",
    "",
    "",
    $deXlf
)
